# Update odds values on row 2 of the "Jogos da Semana" worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value  = 1.75
$ws.Range("H2").Value  = 3.2
$ws.Range("I2").Value  = 5.5
$ws.Range("J2").Value  = 2.5
$ws.Range("K2").Value  = 2
$ws.Range("M2").Value  = 1.1
$ws.Range("N2").Value  = 7
$ws.Range("O2").Value  = 1.5
$ws.Range("P2").Value  = 2.5
$ws.Range("Q2").Value  = 2.5
$ws.Range("R2").Value  = 1.5
$ws.Range("S2").Value  = 1.53
$ws.Range("T2").Value  = 2.38
$ws.Range("U2").Value  = 2.25
$ws.Range("V2").Value  = 1.57
$ws.Range("W2").Value  = 5
$ws.Range("Y2").Value  = 9.5
$ws.Range("Z2").Value  = 13
$ws.Range("AA2").Value = 17
$ws.Range("AB2").Value = 41
$ws.Range("AC2").Value = 6.5
$ws.Range("AE2").Value = 21
$ws.Range("AF2").Value = 81
$ws.Range("AG2").Value = 11
$ws.Range("AJ2").Value = 51
$ws.Range("AO2").Value = 10
$ws.Range("AP2").Value = 26
$ws.Range("AQ2").Value = 34
$ws.Range("AR2").Value = 67
$ws.Range("AS2").Value = 251
$ws.Range("AT2").Value = 2.38
$ws.Range("AU2").Value = 10
$ws.Range("AV2").Value = 81
$ws.Range("AW2").Value = 6.5
